$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.092.10"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "1.907.28"
$ws.Range("E3").Value = "  +5.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.21"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5082"
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.05"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3020"
$ws.Range("E9").Value = "  +8.54%  "
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").Value = "1.907.13"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07317"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6927"
$ws.Range("E14").Value = "  +7.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.62"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.902"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("D17").Value = "30.074.07"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008181"
$ws.Range("E18").Value = "  +11.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +6.23%  "
$ws.Range("D21").Value = "2.153.96"
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.816"
$ws.Range("E23").Value = "  +5.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.730"
$ws.Range("E24").Value = "  +7.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.290"
$ws.Range("E25").Value = "  +4.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.53"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.76"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.05"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("E29").Value = "  +6.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.245"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08806"
$ws.Range("E32").Value = "  +5.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.994"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05049"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7179"
$ws.Range("E36").Value = "  +6.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.691"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.811"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.267"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9644"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01690"
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.128"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4301"
$ws.Range("E43").Value = "  +5.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.43"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.607"
$ws.Range("E46").Value = "  +6.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1278"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05742"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.12"
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.442"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3814"
$ws.Range("E51").Value = "  +5.20%  "
